# "HTML with latest data"
#
# The acharyan captions (column B, rows 3-37 of the "acharyan_captions"
# sheet) had a spelling fix applied throughout: "Svatantra" -> "Swatantra".
# Re-saving the workbook also left the "acharyan_captions" sheet as the
# active/selected tab instead of "Founders_Early_Acharyas".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("acharyan_captions")

for ($row = 3; $row -le 37; $row++) {
    $cell = $ws.Range("B$row")
    $current = $cell.Value2
    if ($current -ne $null -and $current.ToString().Contains("Svatantra")) {
        $cell.Value2 = $current.ToString().Replace("Svatantra", "Swatantra")
    }
}

# Bring the first sheet back into focus (matches the saved workbookView /
# sheetView tabSelected state in the target file).
$ws.Activate()
